$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add a new header label in A1
$ws.Range("A1").Value = "Category"

# Give A1 the same (header) formatting as the other header cells in row 1
$ws.Range("B1").Copy()
$ws.Range("A1").PasteSpecial(-4122)

# The category cells in A2:A46 should no longer carry the header style
$ws.Range("A2:A46").Style = "Normal"
